$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 12: lin2srgb - update Scalpel Type text and mark as a Win ---
$ws.Range("E12").Value = "{'float', 'any', 'int'}"
$ws.Range("F12").Value = "Win"
$ws.Range("F12").Interior.Color = $ws.Range("F11").Interior.Color

# --- Row 13: lin2srgb - update Scalpel Type text and mark as a Win ---
$ws.Range("E13").Value = "float"
$ws.Range("F13").Value = "Win"
$ws.Range("F13").Interior.Color = $ws.Range("F11").Interior.Color

# --- Row 68: Scalpel Wins count goes from 5 to 7 ---
$ws.Range("F68").Value = 7

# --- Insert a new row 69 for "Scalpel Accuracy:" summary, pushing the
#     existing "Accuracy over PyType" row down to row 70 ---
$ws.Rows.Item(69).Insert()
$ws.Range("C69").Value = "Scalpel Accuracy:"
$ws.Range("D69").Value = 100

# Dimension / used range is recalculated automatically by Excel.
